# Update loading_percent values in column F (rows 2-25) on the active sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 16.86991607391245
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("F4").Value = 15.2699775710849
$ws.Range("F5").Value = 15.008197319934
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("F7").Value = 15.26647399323134
$ws.Range("F8").Value = 16.53996406344772
$ws.Range("F9").Value = 19.0027458068253
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("F11").Value = 21.39172256362241
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("F13").Value = 21.60004134736741
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("F15").Value = 21.29868154950794
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("F17").Value = 20.20408069597326
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("F21").Value = 21.46857628470571
$ws.Range("F22").Value = 22.22866616901554
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("F24").Value = 20.22900810905285
$ws.Range("F25").Value = 18.34778573295695
